$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Price (D) and Volume(1h) (E) columns for the crypto rows that changed.
# NumberFormat is set to Text ("@") before assigning so that values such as
# "243.20" or "1.001" keep their exact textual representation (trailing zeros,
# multiple dots, etc.) instead of being auto-converted to numbers by Excel.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.414.85"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  +0.16%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.865.57"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -0.77%  "
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.7064"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -0.76%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "243.20"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  +0.16%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3138"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -1.36%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07852"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "24.45"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -2.55%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.08026"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -3.75%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.884.39"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -1.23%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "93.30"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -1.72%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.6999"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -2.52%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "6.448"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  +1.32%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.000008362"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -2.96%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "29.431.94"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  +0.21%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "252.19"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  +3.72%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "2.126.50"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -0.61%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.12"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -1.60%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.606"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -2.80%  "
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -0.17%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.1558"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -1.08%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.000"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -1.15%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "160.61"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -1.66%  "
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  +0.52%  "
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -0.58%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.320"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -2.63%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.279"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -1.23%  "
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +0.76%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.05306"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.885"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  -3.08%  "
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  -2.56%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.165"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -1.95%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.713"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  +0.91%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01876"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -0.88%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.267.48"
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  +0.14%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.740"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -0.45%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.8996"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -1.04%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "109.07"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -4.41%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "5.977"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -8.10%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "71.33"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -4.51%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.00000000128"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -2.33%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.031.73"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -0.29%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "9.575"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  +0.44%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.787"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -1.23%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.5171"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -1.03%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.4304"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -1.70%  "
